# Add a new "msg_count_twitter_engage" metric block between the existing
# "msg_count_twitter" and "msg_count_facebook" blocks, and refresh the
# Twitter stats with the new numbers (per commit: "new twitter social
# analysis (engage vs museum tweets)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Make room: insert 13 new columns right before the facebook block
#        (old O:AA). This shifts O:AA -> AB:AN, carrying over values,
#        number formats, styles and the O1:AA1 merge intact.
$ws.Range("O1:AA1").EntireColumn.Insert()

# --- 2) Header row 1: merge the freshly inserted O1:AA1 block, match the
#        look of the other two header cells, then label it
$hdr = $ws.Range("O1:AA1")
$hdr.Merge()
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1
$ws.Range("O1").Value = "msg_count_twitter_engage"

# --- 3) Header row 2: stat labels for the new block (same 13 labels as the
#        other two blocks)
$ws.Range("O2").Value = "sum"
$ws.Range("P2").Value = "mean"
$ws.Range("Q2").Value = "std"
$ws.Range("R2").Value = "min"
$ws.Range("S2").Value = "q25"
$ws.Range("T2").Value = "median"
$ws.Range("U2").Value = "q75"
$ws.Range("V2").Value = "max"
$ws.Range("W2").Value = "count"
$ws.Range("X2").Value = "msg_per_mus"
$ws.Range("Y2").Value = "active_mus_n"
$ws.Range("Z2").Value = "active_mus_pc"
$ws.Range("AA2").Value = "active_mus_pc_z"

# --- 4) Refresh the Twitter (B:N) stats with the new figures
$ws.Range("B4").Value = 844126
$ws.Range("C4").Value = 1061.8
$ws.Range("D4").Value = 1994.7
$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 323
$ws.Range("H4").Value = 1373.5
$ws.Range("I4").Value = 19900
$ws.Range("J4").Value = 795
$ws.Range("K4").Value = 1658.4
$ws.Range("L4").Value = 509
$ws.Range("M4").Value = 64
$ws.Range("N4").Value = 0.3

$ws.Range("B5").Value = 1379444
$ws.Range("C5").Value = 575.5
$ws.Range("D5").Value = 2000.7
$ws.Range("E5").Value = 0
$ws.Range("F5").Value = 0
$ws.Range("G5").Value = 43
$ws.Range("H5").Value = 564
$ws.Range("I5").Value = 47580
$ws.Range("J5").Value = 2397
$ws.Range("K5").Value = 981.8
$ws.Range("L5").Value = 1405
$ws.Range("M5").Value = 58.6
$ws.Range("N5").Value = -0.1

$ws.Range("B6").Value = 159975
$ws.Range("C6").Value = 1738.9
$ws.Range("D6").Value = 2237.9
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 36.8
$ws.Range("G6").Value = 960
$ws.Range("H6").Value = 2469.2
$ws.Range("I6").Value = 13060
$ws.Range("J6").Value = 92
$ws.Range("K6").Value = 2253.2
$ws.Range("L6").Value = 71
$ws.Range("M6").Value = 77.2
$ws.Range("N6").Value = 1.1

$ws.Range("B7").Value = 20284
$ws.Range("C7").Value = 368.8
$ws.Range("D7").Value = 1006.2
$ws.Range("E7").Value = 0
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 0
$ws.Range("H7").Value = 107.5
$ws.Range("I7").Value = 4946
$ws.Range("J7").Value = 55
$ws.Range("K7").Value = 922
$ws.Range("L7").Value = 22
$ws.Range("M7").Value = 40
$ws.Range("N7").Value = -1.3

# --- 5) Populate the new "msg_count_twitter_engage" data block (O:AA)
$ws.Range("O4").Value = 315647
$ws.Range("P4").Value = 397
$ws.Range("Q4").Value = 1305.8
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 0
$ws.Range("T4").Value = 40
$ws.Range("U4").Value = 266.5
$ws.Range("V4").Value = 18016
$ws.Range("W4").Value = 795
$ws.Range("X4").Value = 614.1
$ws.Range("Y4").Value = 514
$ws.Range("Z4").Value = 64.7
$ws.Range("AA4").Value = 0.3

$ws.Range("O5").Value = 495606
$ws.Range("P5").Value = 206.8
$ws.Range("Q5").Value = 957.5
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 8
$ws.Range("U5").Value = 103
$ws.Range("V5").Value = 20884
$ws.Range("W5").Value = 2397
$ws.Range("X5").Value = 346.1
$ws.Range("Y5").Value = 1432
$ws.Range("Z5").Value = 59.7
$ws.Range("AA5").Value = 0

$ws.Range("O6").Value = 94956
$ws.Range("P6").Value = 1032.1
$ws.Range("Q6").Value = 5053
$ws.Range("R6").Value = 0
$ws.Range("S6").Value = 5.2
$ws.Range("T6").Value = 167
$ws.Range("U6").Value = 483.5
$ws.Range("V6").Value = 46992
$ws.Range("W6").Value = 92
$ws.Range("X6").Value = 1337.4
$ws.Range("Y6").Value = 71
$ws.Range("Z6").Value = 77.2
$ws.Range("AA6").Value = 1

$ws.Range("O7").Value = 2853
$ws.Range("P7").Value = 51.9
$ws.Range("Q7").Value = 132.5
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 0
$ws.Range("U7").Value = 10
$ws.Range("V7").Value = 648
$ws.Range("W7").Value = 55
$ws.Range("X7").Value = 142.6
$ws.Range("Y7").Value = 20
$ws.Range("Z7").Value = 36.4
$ws.Range("AA7").Value = -1.4
